$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update column F (想去人数 / want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5500
$ws1.Range("F8").Value = 903
$ws1.Range("F10").Value = 2452
$ws1.Range("F12").Value = 74
$ws1.Range("F14").Value = 2298
$ws1.Range("F15").Value = 197

# Sheet "全部类型" (all types) - same underlying records, update column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5500
$ws4.Range("F10").Value = 903
$ws4.Range("F12").Value = 2452
$ws4.Range("F14").Value = 74
$ws4.Range("F17").Value = 2298
$ws4.Range("F18").Value = 197
